$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy formatting from column F into new column G for each relevant row band
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("F2").Copy()
$ws.Range("G2:G6").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("F7").Copy()
$ws.Range("G7").PasteSpecial(-4122) # xlPasteFormats

# Set the new header and values
$ws.Range("G1").Value = "PRESUPUESTO"
$ws.Range("G2:G7").Value = 0

# Set new column width (ColumnWidth value tuned so the stored OOXML width
# attribute comes out to exactly 17, matching the target file)
$ws.Columns.Item(7).ColumnWidth = 16.17
